$wb = $excel.ActiveWorkbook

# --- Sheets ---
$wsAttendance = $wb.Worksheets.Item(1)   # "Điểm danh"
$wsLeave      = $wb.Worksheets.Item(2)   # "Nghỉ phép"

# --- Rename the approval-status headers on the "Nghỉ phép" sheet ---
$wsLeave.Range("G1").Value = "Ngày Duyệt/Từ chối Lần đầu"
$wsLeave.Range("H1").Value = "Trạng thái lần Lần đầu"
$wsLeave.Range("I1").Value = "Ngày Duyệt/Từ chối Lần cuối"
$wsLeave.Range("J1").Value = "Trạng thái Lần cuối"

# --- Widen a couple of columns on the "Nghỉ phép" sheet ---
$wsLeave.Columns.Item(1).ColumnWidth = 6
$wsLeave.Columns.Item(10).ColumnWidth = 15.571428571428571

# --- Make "Nghỉ phép" the active/selected sheet (was "Điểm danh") ---
$wsLeave.Activate() | Out-Null
$wsLeave.Range("A2").Select() | Out-Null
